# Apply "add links to chapters 34.1, 31.1, 30.2 to labs 14, 06, 01 respectively"
# on sheet "Лист1" of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Chapter 31.1 is renumbered to 31 (its row now links to lab 06 / column E).
$ws.Range("A35").Value = "31"
$ws.Range("E35").Value = 1

# Chapter 34 is renumbered to 34.1 (its row now links to lab 14 / column L).
$ws.Range("A38").Value = "34.1"
$ws.Range("L38").Value = 1

# Chapter 30 (row 34) now also links to lab 01 / column B.
$ws.Range("B34").Value = 1

# Chapter 33 (row 37) now also links to lab 01 / column B.
$ws.Range("B37").Value = 1

# Match the cursor position left behind by the author's edit.
$ws.Range("P32").Select()
